$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3453
$ws.Range("F5").Value = 8081
$ws.Range("F7").Value = 60
$ws.Range("F8").Value = 2065
$ws.Range("F11").Value = 531
$ws.Range("F12").Value = 433
$ws.Range("F14").Value = 1043
$ws.Range("F16").Value = 144
$ws.Range("F17").Value = 1119
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 411
$ws.Range("F23").Value = 4430
$ws.Range("F24").Value = 85
$ws.Range("F25").Value = 46901
$ws.Range("F26").Value = 3898
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 980
$ws.Range("F29").Value = 699
$ws.Range("F30").Value = 4
$ws.Range("F32").Value = 815
$ws.Range("F35").Value = 181
$ws.Range("F38").Value = 820
$ws.Range("F39").Value = 923
$ws.Range("F40").Value = 109
$ws.Range("F43").Value = 666
$ws.Range("F44").Value = 80
$ws.Range("F45").Value = 62
$ws.Range("F46").Value = 13
$ws.Range("F47").Value = 2440

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 7249

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2189
$ws.Range("F5").Value = 1463
$ws.Range("F10").Value = 1474
$ws.Range("F12").Value = 52

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3453
$ws.Range("F3").Value = 2189
$ws.Range("F4").Value = 8081
$ws.Range("F5").Value = 1463
$ws.Range("F8").Value = 52
$ws.Range("F9").Value = 60
$ws.Range("F11").Value = 531
$ws.Range("F12").Value = 1043
$ws.Range("F15").Value = 1119
$ws.Range("F17").Value = 55
$ws.Range("F18").Value = 4430
$ws.Range("F19").Value = 85
$ws.Range("F23").Value = 3898
$ws.Range("F25").Value = 980
$ws.Range("F26").Value = 699
$ws.Range("F28").Value = 815
$ws.Range("F32").Value = 181
$ws.Range("F36").Value = 923
$ws.Range("F38").Value = 109
$ws.Range("F41").Value = 666
$ws.Range("F43").Value = 80
$ws.Range("F45").Value = 62
$ws.Range("F46").Value = 13
$ws.Range("F48").Value = 2440

